$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: extend header with two new columns (P1, Q1) ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the existing header formatting (bold, centered, top-aligned, bordered)
# from O1 onto the two new header cells so they share the same style index.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)   # xlPasteFormats

# --- Rows 2-25: update data cells ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
